$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Scalar value corrections (floating point precision fixes) ---
$ws.Range("T57").Value = 2516232.96
$ws.Range("AJ57").Value = 4558894.592
$ws.Range("H59").Value = 1825236.096
$ws.Range("AB59").Value = 1915573.12
$ws.Range("AF59").Value = 2934638.08
$ws.Range("AJ59").Value = 3678493.952
$ws.Range("AF60").Value = -1949687.04
$ws.Range("AJ60").Value = -2333126.144
$ws.Range("X61").Value = 1001162.048
$ws.Range("AB61").Value = 833234.112
$ws.Range("AF61").Value = 984950.912
$ws.Range("AJ61").Value = 1345367.936
$ws.Range("AB62").Value = -73811
$ws.Range("AF63").Value = -140338
$ws.Range("X66").Value = -66438.008
$ws.Range("AF66").Value = -202331.984
$ws.Range("AJ66").Value = -138730.032
$ws.Range("X67").Value = -417337.984
$ws.Range("AF67").Value = 29634
$ws.Range("P69").Value = 14953.992
$ws.Range("L71").Value = 2854
$ws.Range("L72").Value = 133200.024
$ws.Range("X72").Value = -1378766.08
$ws.Range("AJ72").Value = 19090.984
$ws.Range("P74").Value = 64757.992
$ws.Range("X74").Value = 164211.024
$ws.Range("AB74").Value = 474533.952
$ws.Range("AF74").Value = 406356.928
$ws.Range("AJ74").Value = 965980.032
$ws.Range("AB75").Value = -68715.992
$ws.Range("AF75").Value = 54785.016
$ws.Range("AJ75").Value = 22728.992
$ws.Range("AB79").Value = -75982.008
$ws.Range("L80").Value = 61321.968
$ws.Range("T80").Value = 118154.984
$ws.Range("AF80").Value = 397159.968
$ws.Range("AJ80").Value = 820111.936

# --- Row 64: clear placeholder zero values (now blank cells) ---
$ws.Range("C64").ClearContents()
$ws.Range("E64:AK64").ClearContents()

# --- Row 79: clear placeholder zero values (now blank cells) ---
$ws.Range("C79").ClearContents()
$ws.Range("E79:X79").ClearContents()
